# The deck's single editable theme (backed by ppt/theme/theme1.xml,
# linked from the slide master that every slide/layout uses) is switched
# from the custom "Integral" palette to the stock "Office" palette.
#
# PowerPoint's Design/ColorScheme COM surface only lets us poke individual
# RGB values (ColorScheme.Colors(i).RGB, in the usual COM BGR-packed-int
# form), so we replay the 12 theme colour slots in order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3, 8 accent4,
#   9 accent5, 10 accent6, 11 hlink, 12 folHlink
$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $cs.Colors($i).RGB = $officeThemeColors[$i - 1]
}
